$wb = $excel.ActiveWorkbook

# --- Sheet "Persone-Indisp": append new person rows (A6:A13) ---
$wsPersone = $wb.Worksheets.Item("Persone-Indisp")

$wsPersone.Cells.Item(6, 1).Value = "MAR"
$wsPersone.Cells.Item(7, 1).Value = "MIG"
$wsPersone.Cells.Item(8, 1).Value = "FAN"
$wsPersone.Cells.Item(9, 1).Value = "LEG"
$wsPersone.Cells.Item(10, 1).Value = "SAR"
$wsPersone.Cells.Item(11, 1).Value = "BAT"
$wsPersone.Cells.Item(13, 1).Value = "BOM"
$wsPersone.Cells.Item(12, 1).Value = "AIN"

$wsPersone.Range("K14").Select()

# --- Sheet "Turni Fissi": fill weekday shifts with "BAI" for weeks 2-4 ---
$wsTurni = $wb.Worksheets.Item("Turni Fissi")

$rows = @(13, 14, 15, 16, 17, 20, 21, 22, 23, 24, 27, 28, 29, 30, 31)
foreach ($r in $rows) {
    $wsTurni.Cells.Item($r, 2).Value = "BAI"
}

$wsTurni.Range("D6").Select()

$wsPersone.Activate()
